{"js": "const replacements = [\n  [\"2024-10-16 Wednesday\", \"2024-10-17 Thursday\"],\n  [\"217\u00d77=1519\", \"496\u00d74=1984\"],\n  [\"682\u00d78=5456\", \"736\u00d74=2944\"],\n  [\"239\u00d78=1912\", \"573\u00d73=1719\"],\n  [\"654\u00d75=3270\", \"163\u00d78=1304\"],\n  [\"460\u00d77=3220\", \"748\u00d73=2244\"],\n  [\"279\u00d74=1116\", \"434\u00d73=1302\"],\n  [\"210\u00d76=1260\", \"859\u00d76=5154\"],\n  [\"384\u00d73=1152\", \"644\u00d78=5152\"],\n  [\"367\u00d76=2202\", \"677\u00d72=1354\"],\n  [\"366\u00d74=1464\", \"784\u00d76=4704\"],\n  [\"157\u00d77=1099\", \"408\u00d73=1224\"],\n  [\"674\u00d76=4044\", \"342\u00d72=684\"],\n  [\"386\u00d79=3474\", \"516\u00d72=1032\"],\n  [\"351\u00d75=1755\", \"233\u00d79=2097\"],\n  [\"761\u00d73=2283\", \"833\u00d74=3332\"],\n  [\"740\u00d76=4440\", \"521\u00d75=2605\"],\n  [\"159\u00d76=954\", \"419\u00d74=1676\"],\n  [\"842\u00d75=4210\", \"237\u00d78=1896\"],\n  [\"818\u00d73=2454\", \"561\u00d79=5049\"],\n  [\"968\u00d75=4840\", \"185\u00d76=1110\"],\n  [\"879\u00d72=1758\", \"506\u00d75=2530\"],\n  [\"466\u00d72=932\", \"632\u00d76=3792\"],\n  [\"214\u00d76=1284\", \"145\u00d73=435\"],\n  [\"180\u00d73=540\", \"260\u00d73=780\"],\n  [\"205\u00d77=1435\", \"144\u00d78=1152\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-10-16 Wednesday\", \"2024-10-17 Thursday\"),\n    @(\"217\u00d77=1519\", \"496\u00d74=1984\"),\n    @(\"682\u00d78=5456\", \"736\u00d74=2944\"),\n    @(\"239\u00d78=1912\", \"573\u00d73=1719\"),\n    @(\"654\u00d75=3270\", \"163\u00d78=1304\"),\n    @(\"460\u00d77=3220\", \"748\u00d73=2244\"),\n    @(\"279\u00d74=1116\", \"434\u00d73=1302\"),\n    @(\"210\u00d76=1260\", \"859\u00d76=5154\"),\n    @(\"384\u00d73=1152\", \"644\u00d78=5152\"),\n    @(\"367\u00d76=2202\", \"677\u00d72=1354\"),\n    @(\"366\u00d74=1464\", \"784\u00d76=4704\"),\n    @(\"157\u00d77=1099\", \"408\u00d73=1224\"),\n    @(\"674\u00d76=4044\", \"342\u00d72=684\"),\n    @(\"386\u00d79=3474\", \"516\u00d72=1032\"),\n    @(\"351\u00d75=1755\", \"233\u00d79=2097\"),\n    @(\"761\u00d73=2283\", \"833\u00d74=3332\"),\n    @(\"740\u00d76=4440\", \"521\u00d75=2605\"),\n    @(\"159\u00d76=954\", \"419\u00d74=1676\"),\n    @(\"842\u00d75=4210\", \"237\u00d78=1896\"),\n    @(\"818\u00d73=2454\", \"561\u00d79=5049\"),\n    @(\"968\u00d75=4840\", \"185\u00d76=1110\"),\n    @(\"879\u00d72=1758\", \"506\u00d75=2530\"),\n    @(\"466\u00d72=932\", \"632\u00d76=3792\"),\n    @(\"214\u00d76=1284\", \"145\u00d73=435\"),\n    @(\"180\u00d73=540\", \"260\u00d73=780\"),\n    @(\"205\u00d77=1435\", \"144\u00d78=1152\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
